$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Info")

# Keep "Info" as the active/selected sheet tab, matching the original file.
$ws.Activate()

# Insert a new column before column A, shifting the existing data
# (columns A-C) right to columns B-D.
$ws.Columns.Item(1).Insert()

# Populate the newly inserted column A with a header label.
$ws.Range("A1").Value = "Metadata"

# Move the active selection to A2 (re-envisioned loading sheet tab logic).
$ws.Range("A2").Select()

# Preserve intended tab-scroll-bar ratio from the source edit.
$wb.Windows.Item(1).TabRatio = 0.494

# Rename the duplicated "Excel Built-in" cell style to match upstream.
foreach ($s in $wb.Styles) {
    if ($s.Name -eq "Excel Built-in Excel Built-in Normal 2") {
        $s.Name = "Excel Built-in Excel Built-in Excel Built-in Normal 2"
    }
}
